$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E5 value (was "PUR", now "PUR : 35 | SUPP : 65")
$ws.Range("E5").Value = "PUR : 35 | SUPP : 65"

# Clear the now-unused column K cells (K3, K4, K5) entirely
$ws.Range("K3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()

# Move the active selection to E6 (was K6)
$ws.Range("E6").Select()
